$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write values in the same order the original workbook's shared-string
# table was built (so new unique strings land at matching indices).
$ws.Range("A1").Value = "from"
$ws.Range("B1").Value = "where"
$ws.Range("C1").Value = "Mobile_no"
$ws.Range("D1").Value = "Mobile_no"
$ws.Range("E1").Value = "error_message"
$ws.Range("F1").Value = "from"
$ws.Range("G1").Value = "to"
$ws.Range("H1").Value = "adults"
$ws.Range("I1").Value = "children"
$ws.Range("J1").Value = "infants"
$ws.Range("K1").Value = "travel_class"

$ws.Range("A2").Value = "Bengaluru"
$ws.Range("B2").Value = "Manali"
$ws.Range("C2").Value = 8015993932
$ws.Range("D2").Value = 12345
$ws.Range("F2").Value = "Chennai"
$ws.Range("G2").Value = "Mumbai"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "Economy"
$ws.Range("E2").Value = "Please enter a valid number"

# Column width for C (matches the saved width="11" after Excel's
# char-width rounding)
$ws.Columns.Item(3).ColumnWidth = 10.1

# Selection
$ws.Range("L1:L4").Select()
